$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("product_import_template")

# Update the Product Code value in A2
$ws.Range("A2").Value = 102048

# Update Category (F2) text from MEDICINE to MEDICINES
$ws.Range("F2").Value = "MEDICINES"

# Unhide column A
$ws.Columns.Item(1).Hidden = $false

# Update the sheet view: remove frozen/scrolled topLeftCell (scroll back to A1) and
# change the active selection to A2
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A2").Select()
